$d = $word.ActiveDocument

# 1. "Pendekatan diatas memungkinkan" -> "Pendekatan di atas memungkinkan"
#    (splits the word "diatas" into "di atas")
$d.Content.Find.Execute(
    "Pendekatan diatas memungkinkan", $false, $false, $false, $false, $false,
    $true, 1, $false, "Pendekatan di atas memungkinkan", 2) | Out-Null

# 2. Replace the graphics-library list.
#    " yang meliputi Windows Header,  GLUT, dan ImGUI. "
#      -> " yang meliputi OpenGL, GLFW, GLSL , dan ImGUI. "
#    Done as two Find/Replace calls that stop right at the "_GoBack" bookmark
#    so the bookmark (which sits between "meliputi " and "Windows Header")
#    stays in place between "OpenGL," and " GLFW, GLSL ,".
$d.Content.Find.Execute(
    "yang meliputi ", $false, $false, $false, $false, $false,
    $true, 1, $false, "yang meliputi OpenGL,", 2) | Out-Null

$d.Content.Find.Execute(
    "Windows Header,  GLUT, dan ImGUI. ", $false, $false, $false, $false, $false,
    $true, 1, $false, " GLFW, GLSL , dan ImGUI. ", 2) | Out-Null

# 3. "...sebagai berikut :" -> "...sebagai berikut:" (drop the space before the colon)
$d.Content.Find.Execute(
    "Adapun sistematika penulisan yang digunakan dalam penulisan ini adalah sebagai berikut :",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "Adapun sistematika penulisan yang digunakan dalam penulisan ini adalah sebagai berikut:", 2) | Out-Null
